$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Angle" header column (AG1), copying the header style from AF1
$ws.Range("AF1:AF1").Copy($ws.Range("AG1"))
$ws.Range("AG1").Value = "Angle"

# Row 2
$ws.Range("K2").Value = "['#5c9909', '#3e7204', '#7cac3b']"
$ws.Range("L2").Value = 91.85866616780461
$ws.Range("M2").Value = 153.3836872135442
$ws.Range("N2").Value = 8.64465438219056
$ws.Range("P2").Value = 0.5231236527532824
$ws.Range("Q2").Value = 62.47787610619413
$ws.Range("R2").Value = 114.3943362831863
$ws.Range("S2").Value = 3.942300884955552
$ws.Range("U2").Value = 0.2776308054085832
$ws.Range("V2").Value = 124.0746745271427
$ws.Range("W2").Value = 171.8636698599856
$ws.Range("X2").Value = 59.33284205354913
$ws.Range("Z2").Value = 0.1992455418381344
$ws.Range("AA2").Value = 425
$ws.Range("AB2").Value = 312
$ws.Range("AC2").Value = 29075
$ws.Range("AD2").Value = 1.362179487179487
$ws.Range("AE2").Value = 0.2192684766214178
$ws.Range("AF2").Value = 0.7917273680340926
$ws.Range("AG2").Value = 28.65376091003418

# Row 3
$ws.Range("K3").Value = "['#1f2b0d', '#3e4d1f', '#ab90d1']"
$ws.Range("L3").Value = 31.09900725042031
$ws.Range("M3").Value = 42.56052426101493
$ws.Range("N3").Value = 13.09405465699855
$ws.Range("P3").Value = 0.5643728911355569
$ws.Range("Q3").Value = 62.094764198783
$ws.Range("R3").Value = 76.7427527045309
$ws.Range("S3").Value = 31.27389705882275
$ws.Range("U3").Value = 0.2947204696260803
$ws.Range("V3").Value = 171.1346359266632
$ws.Range("W3").Value = 144.4651949439139
$ws.Range("X3").Value = 209.2709631475868
$ws.Range("Z3").Value = 0.1409066392383628
$ws.Range("AA3").Value = 499
$ws.Range("AB3").Value = 332
$ws.Range("AC3").Value = 160978
$ws.Range("AD3").Value = 1.503012048192771
$ws.Range("AE3").Value = 0.9716903686891856
$ws.Range("AF3").Value = 0.9716903686891856
$ws.Range("AG3").Value = 81.53978729248047

# Row 4
$ws.Range("K4").Value = "['#152412', '#2c5221', '#497c3e']"
$ws.Range("L4").Value = 21.4819562330384
$ws.Range("M4").Value = 36.44638835960434
$ws.Range("N4").Value = 17.53143945955303
$ws.Range("P4").Value = 0.4417112299465241
$ws.Range("Q4").Value = 43.52746689448949
$ws.Range("R4").Value = 81.56744980777452
$ws.Range("S4").Value = 32.53208030756082
$ws.Range("T4").Value = "darkslategray"
$ws.Range("U4").Value = 0.297886427298192
$ws.Range("V4").Value = 73.30060511419153
$ws.Range("W4").Value = 124.1522545383554
$ws.Range("X4").Value = 61.59037673238332
$ws.Range("Z4").Value = 0.260402342755284
$ws.Range("AA4").Value = 499
$ws.Range("AB4").Value = 332
$ws.Range("AC4").Value = 38861
$ws.Range("AD4").Value = 1.503012048192771
$ws.Range("AE4").Value = 0.2345715527440423
$ws.Range("AF4").Value = 0.6625690513537476
$ws.Range("AG4").Value = 155.7671661376953

# Row 5
$ws.Range("K5").Value = "['#787839', '#3f3819', '#b9b683']"
$ws.Range("L5").Value = 119.6535114850098
$ws.Range("M5").Value = 120.3929234195226
$ws.Range("N5").Value = 57.46525629922908
$ws.Range("P5").Value = 0.4678149625262509
$ws.Range("Q5").Value = 62.69727001884973
$ws.Range("R5").Value = 56.23213959992711
$ws.Range("S5").Value = 24.59135404633094
$ws.Range("U5").Value = 0.3162437624029439
$ws.Range("V5").Value = 184.7987864727414
$ws.Range("W5").Value = 181.5444811278665
$ws.Range("X5").Value = 130.5095922191479
$ws.Range("Y5").Value = "darkkhaki"
$ws.Range("Z5").Value = 0.2159412750708052
$ws.Range("AA5").Value = 387
$ws.Range("AC5").Value = 53230
$ws.Range("AD5").Value = 1.6125
$ws.Range("AE5").Value = 0.573105081826012
$ws.Range("AF5").Value = 0.7515973031169473
$ws.Range("AG5").Value = 104.5516662597656

